# Update the logged "Lopetus" (end) time for the entry on row 15 - this
# pushes the row's Tunnit (hours) and the running Total Tunnit formulas
# (D15 and F15:F26, which all chain off of it) to new recalculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

$ws.Range("C15").Value = 1.0208333333333333

# Move the active selection to D11, matching the sheet's last-saved cursor
# position.
$ws.Range("D11").Select() | Out-Null
